# Apply dSF (column F) corrections from repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -5
    "F4"  = 3
    "F9"  = -3
    "F13" = -1
    "F21" = 4
    "F23" = 4
    "F26" = -5
    "F27" = -1
    "F29" = -3
    "F32" = 3
    "F35" = 3
    "F37" = -3
    "F40" = -1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
